# Append/refresh scrape: new fetch timestamp, drop a filtered-out listing (old row 24),
# shift the remaining listings up by one row, and drop the now-empty last row.
#
# Before -> After mapping (sheet "ランサーズ"):
#   rows 2-23 : same listing, timestamp refreshed
#   row 24    : takes the listing that used to be on row 25 (AWS infra engineer)
#   row 25    : takes the listing that used to be on row 26 (Excel macro)
#   row 26    : removed (the sheet shrinks from A1:H26 to A1:H25)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-21 01:19:03"

# 1) Refresh the "取得日時" (fetched at) timestamp for every existing data row (2..26).
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# 2) Old row 24 (a filtered-out/spam listing) is replaced by what used to be row 25.
$ws.Cells.Item(24, 2).Value = "【継続案件あり】AWSに精通しているインフラエンジニアを募集します"
$ws.Cells.Item(24, 6).Value = "https://www.lancers.jp/work/detail/5416510"

# 3) Old row 25 is replaced by what used to be row 26.
$ws.Cells.Item(25, 2).Value = "【急募】エクセルマクロの組み方を教えてください!"
$ws.Cells.Item(25, 6).Value = "https://www.lancers.jp/work/detail/5416433"

# 4) The old row 26 no longer exists after the shift - remove it (also fixes dimension/UsedRange).
$ws.Rows.Item(26).Delete()

# 5) Rebuild the hyperlinks collection so F2:F25 point at the correct (possibly shifted) URLs,
#    and the stale hyperlink that used to live on F26 is gone.
$ws.Hyperlinks.Delete()
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
}
